$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (row 26) and the "SC 92" row (originally row 28,
# becomes row 27 after the first deletion shifts rows up). All subsequent
# rows shift up by one each time, ending up two rows higher overall.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# After the shift, a handful of column D ("E" on-sheet) values were
# re-imputed differently than before, so fix those three cells up.
$ws.Range("E27").Value = -10      # SC 101 row: was blank, now has a value
$ws.Range("E29").ClearContents()  # SC 119 row: was -6.8, now blank
$ws.Range("E32").ClearContents()  # SC 193 row: was -6.4, now blank
